# 21_day_menu_database_wk3_tuesday_with_usda.xlsx
# "simplified plus add sodium and fiber"
#
# The underlying data (columns A:T, rows 1:18 - including the already
# present Fiber (g)/Sodium (mg) columns N & P) is unchanged by this
# revision; the meaningful, user-visible edits captured by the diff are:
#   1. The worksheet tab is renamed from "Sheet1" to "3_Tues".
#   2. The active/selected cell is moved to G1 (the Item Name column).
#   3. Columns A (Date) and G (Item Name) are resized to fit their
#      contents ("bestFit"/"customWidth" columns added to the sheet).
#
# (Everything else in the raw OOXML diff - refreshed fileVersion/rupBuild
# markers, mc:AlternateContent/xr:revisionPtr blocks, bookViews window
# geometry, styles.xml/theme1.xml namespace & naming churn, x14ac:dyDescent
# decoration, and a handful of cells whose literal float text grew extra
# trailing digits such as 141.8 -> 141.80000000000001 - is the same IEEE
# value re-serialized by a newer Excel build on save, not a deliberate
# content edit, so there is nothing for a COM script to reproduce there.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "3_Tues"

# 2. Select cell G1.
$ws.Range("G1").Select()

# 3. Auto-fit column widths for columns A and G.
$ws.Columns("A").AutoFit()
$ws.Columns("G").AutoFit()
